function Get-ParagraphByText($doc, $needle) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

function Replace-ParagraphXml($doc, $para, $xmlBody) {
    # Replace the full paragraph (including its trailing paragraph mark) with
    # the supplied WordprocessingML body fragment (one or more <w:p> elements).
    $full = $doc.Range($para.Range.Start, $para.Range.End)
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
           '<w:body>' + $xmlBody + '</w:body>' +
           '</w:document>' +
           '</pkg:xmlData>' +
           '</pkg:part>' +
           '</pkg:package>'
    $full.InsertXML($pkg)
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove <w:lastRenderedPageBreak/> from the run that holds
#    "SEÇÃO A – REGIME GERAL (REGRA PADRÃO)"
# ---------------------------------------------------------------------------
$p1 = Get-ParagraphByText $d "SEÇÃO A – REGIME GERAL (REGRA PADRÃO)"
$body1 = '<w:p w14:paraId="3771002E" w14:textId="64AD8C7E" w:rsidR="003C6448" w:rsidRPr="00771E6B" w:rsidRDefault="003C6448" w:rsidP="00771E6B">' +
           '<w:pPr>' +
             '<w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>' +
             '<w:outlineLvl w:val="1"/>' +
             '<w:rPr>' +
               '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
               '<w:b/><w:bCs/><w:kern w:val="0"/><w:lang w:val="en-BR"/><w14:ligatures w14:val="none"/>' +
             '</w:rPr>' +
           '</w:pPr>' +
           '<w:r w:rsidRPr="00771E6B">' +
             '<w:rPr>' +
               '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
               '<w:b/><w:bCs/><w:kern w:val="0"/><w:lang w:val="en-BR"/><w14:ligatures w14:val="none"/>' +
             '</w:rPr>' +
             '<w:t>SEÇÃO A – REGIME GERAL (REGRA PADRÃO)</w:t>' +
           '</w:r>' +
         '</w:p>'
Replace-ParagraphXml $d $p1 $body1

# ---------------------------------------------------------------------------
# 2) Add <w:lastRenderedPageBreak/> to the run that holds
#    "(Aplicável quando houver propriedade formal comprovada ou autorização
#    regular dos titulares)"
# ---------------------------------------------------------------------------
$p2 = Get-ParagraphByText $d "(Aplicável quando houver propriedade formal comprovada ou autorização regular dos titulares)"
$body2 = '<w:p w14:paraId="50BA213E" w14:textId="245F0718" w:rsidR="003C6448" w:rsidRPr="00771E6B" w:rsidRDefault="003C6448" w:rsidP="003C6448">' +
           '<w:pPr>' +
             '<w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>' +
             '<w:rPr>' +
               '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
               '<w:kern w:val="0"/><w:lang w:val="en-BR"/><w14:ligatures w14:val="none"/>' +
             '</w:rPr>' +
           '</w:pPr>' +
           '<w:r w:rsidRPr="00771E6B">' +
             '<w:rPr>' +
               '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
               '<w:i/><w:iCs/><w:kern w:val="0"/><w:lang w:val="en-BR"/><w14:ligatures w14:val="none"/>' +
             '</w:rPr>' +
             '<w:lastRenderedPageBreak/>' +
             '<w:t>(Aplicável quando houver propriedade formal comprovada ou autorização regular dos titulares)</w:t>' +
           '</w:r>' +
         '</w:p>'
Replace-ParagraphXml $d $p2 $body2

# ---------------------------------------------------------------------------
# 3) Insert two new empty (bold, outlineLvl=2) paragraphs between the empty
#    "spacing after=0" paragraph and the "6A. VIGÊNCIA" heading paragraph.
# ---------------------------------------------------------------------------
$p3 = Get-ParagraphByText $d "6A. VIGÊNCIA"
$body3 = '<w:p>' +
           '<w:pPr>' +
             '<w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>' +
             '<w:outlineLvl w:val="2"/>' +
             '<w:rPr>' +
               '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
               '<w:b/><w:bCs/><w:kern w:val="0"/><w:lang w:val="en-BR"/><w14:ligatures w14:val="none"/>' +
             '</w:rPr>' +
           '</w:pPr>' +
         '</w:p>' +
         '<w:p>' +
           '<w:pPr>' +
             '<w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>' +
             '<w:outlineLvl w:val="2"/>' +
             '<w:rPr>' +
               '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
               '<w:b/><w:bCs/><w:kern w:val="0"/><w:lang w:val="en-BR"/><w14:ligatures w14:val="none"/>' +
             '</w:rPr>' +
           '</w:pPr>' +
         '</w:p>' +
         '<w:p w14:paraId="673C1C39" w14:textId="10A41830" w:rsidR="003C6448" w:rsidRPr="00B34707" w:rsidRDefault="003C6448" w:rsidP="00B34707">' +
           '<w:pPr>' +
             '<w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>' +
             '<w:outlineLvl w:val="2"/>' +
             '<w:rPr>' +
               '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
               '<w:b/><w:bCs/><w:kern w:val="0"/><w:lang w:val="en-BR"/><w14:ligatures w14:val="none"/>' +
             '</w:rPr>' +
           '</w:pPr>' +
           '<w:r w:rsidRPr="00771E6B">' +
             '<w:rPr>' +
               '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
               '<w:b/><w:bCs/><w:kern w:val="0"/><w:lang w:val="en-BR"/><w14:ligatures w14:val="none"/>' +
             '</w:rPr>' +
             '<w:lastRenderedPageBreak/>' +
             '<w:t>6A. VIGÊNCIA</w:t>' +
           '</w:r>' +
         '</w:p>'
Replace-ParagraphXml $d $p3 $body3

# ---------------------------------------------------------------------------
# 4) Split the "6B. CIÊNCIA REGULATÓRIA" heading paragraph into an empty
#    paragraph followed by a new paragraph holding the (now page-break-free)
#    heading run.
# ---------------------------------------------------------------------------
$p4 = Get-ParagraphByText $d "6B. CIÊNCIA REGULATÓRIA"
$body4 = '<w:p w14:paraId="44BB7572" w14:textId="21E4AD6A" w:rsidR="003C6448" w:rsidRPr="00B34707" w:rsidRDefault="003C6448" w:rsidP="00B34707">' +
           '<w:pPr>' +
             '<w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>' +
             '<w:outlineLvl w:val="2"/>' +
             '<w:rPr>' +
               '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
               '<w:b/><w:bCs/><w:kern w:val="0"/><w:lang w:val="en-BR"/><w14:ligatures w14:val="none"/>' +
             '</w:rPr>' +
           '</w:pPr>' +
         '</w:p>' +
         '<w:p>' +
           '<w:pPr>' +
             '<w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>' +
             '<w:outlineLvl w:val="2"/>' +
             '<w:rPr>' +
               '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
               '<w:b/><w:bCs/><w:kern w:val="0"/><w:lang w:val="en-BR"/><w14:ligatures w14:val="none"/>' +
             '</w:rPr>' +
           '</w:pPr>' +
           '<w:r>' +
             '<w:rPr>' +
               '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
               '<w:b/><w:bCs/><w:kern w:val="0"/><w:lang w:val="en-BR"/><w14:ligatures w14:val="none"/>' +
             '</w:rPr>' +
             '<w:t>6B. CIÊNCIA REGULATÓRIA</w:t>' +
           '</w:r>' +
         '</w:p>'
Replace-ParagraphXml $d $p4 $body4

# ---------------------------------------------------------------------------
# 5) Split the "9. ASSINATURAS" heading paragraph into an empty paragraph
#    followed by a new paragraph holding the (now page-break-free) heading
#    run.
# ---------------------------------------------------------------------------
$p5 = Get-ParagraphByText $d "9. ASSINATURAS"
$body5 = '<w:p w14:paraId="0E314C4F" w14:textId="74413CBD" w:rsidR="003C6448" w:rsidRPr="00B34707" w:rsidRDefault="003C6448" w:rsidP="00B34707">' +
           '<w:pPr>' +
             '<w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>' +
             '<w:outlineLvl w:val="1"/>' +
             '<w:rPr>' +
               '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
               '<w:b/><w:bCs/><w:kern w:val="0"/><w:lang w:val="en-BR"/><w14:ligatures w14:val="none"/>' +
             '</w:rPr>' +
           '</w:pPr>' +
         '</w:p>' +
         '<w:p>' +
           '<w:pPr>' +
             '<w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>' +
             '<w:outlineLvl w:val="1"/>' +
             '<w:rPr>' +
               '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
               '<w:b/><w:bCs/><w:kern w:val="0"/><w:lang w:val="en-BR"/><w14:ligatures w14:val="none"/>' +
             '</w:rPr>' +
           '</w:pPr>' +
           '<w:r>' +
             '<w:rPr>' +
               '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
               '<w:b/><w:bCs/><w:kern w:val="0"/><w:lang w:val="en-BR"/><w14:ligatures w14:val="none"/>' +
             '</w:rPr>' +
             '<w:t>9. ASSINATURAS</w:t>' +
           '</w:r>' +
         '</w:p>'
Replace-ParagraphXml $d $p5 $body5

Write-Host "Applied all 5 changes."
